# ---------------------------------------------------------------------------
# Integrated Data Provider with the framework for parameterized testing
#
#  - Renames the original "Runner" sheet to "Tests"
#  - Fixes C3 on "Tests" from "no" to "yes"
#  - Bolds + centers the header row on "Tests" and updates the selection
#  - Adds a new "DataProviderTests" sheet with username/password data-driven
#    rows, bold+centered header, and centered data cells
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- "Runner" -> "Tests" -----------------------------------------------
$wsTests = $wb.Worksheets.Item(1)
$wsTests.Name = "Tests"

# Correct the InValid-credentials row: execute flag "no" -> "yes"
$wsTests.Range("C3").Value = "yes"

# Header row: bold + centered (matches new cellXfs entry fontId=1)
$hdrTests = $wsTests.Range("A1:E1")
$hdrTests.Font.Bold = $true
$hdrTests.HorizontalAlignment = -4108  # xlCenter
$hdrTests.VerticalAlignment = -4108    # xlCenter

# New selection / view state left behind on the Tests sheet
$wsTests.Range("B13").Select() | Out-Null

# --- New "DataProviderTests" sheet --------------------------------------
$wsData = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsData.Name = "DataProviderTests"

$wsData.Range("A1").Value = "testname"
$wsData.Range("B1").Value = "execute"
$wsData.Range("C1").Value = "username"
$wsData.Range("D1").Value = "password"

# Reuse the already-built bold+centered header style (avoids generating
# redundant intermediate cellXfs entries that a fresh Font/Alignment
# mutation sequence would leave behind). Copy from a single cell so the
# paste doesn't spill past the 4-column target range.
$wsTests.Range("A1").Copy() | Out-Null
$wsData.Range("A1:D1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

$rows = @(
    @("loginTestWithValidCredentials",   "yes", "Admin",   "admin123"),
    @("loginTestWithValidCredentials",   "no",  "adminnn", "admin123"),
    @("loginTestWithInValidCredentials", "yes", "Admin",   "admin12345"),
    @("loginTestWithValidCredentials",   "yes", "Admin",   "admin123"),
    @("loginTestWithValidCredentials",   "no",  "adminnn", "admin123"),
    @("loginTestWithInValidCredentials", "yes", "Admin",   "admin12345")
)

$r = 2
foreach ($row in $rows) {
    $wsData.Cells.Item($r, 1).Value = $row[0]
    $wsData.Cells.Item($r, 2).Value = $row[1]
    $wsData.Cells.Item($r, 3).Value = $row[2]
    $wsData.Cells.Item($r, 4).Value = $row[3]
    $r++
}

# Data rows reuse the plain centered/middle style (cellXfs index 1) that
# the original sheet's body rows already carry.
$wsTests.Range("A2").Copy() | Out-Null
$wsData.Range("A2:D7").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# Column widths approximating the authored sheet (best-fit sizing isn't
# reproducible bit-for-bit through COM, so set close numeric widths).
$wsData.Columns.Item(1).ColumnWidth = 29.5703125
$wsData.Columns.Item(2).ColumnWidth = 7.140625
$wsData.Columns.Item(3).ColumnWidth = 8.85546875
$wsData.Columns.Item(4).ColumnWidth = 10.5703125

$wsData.Range("D14").Select() | Out-Null

$wsData.Application.ActiveWindow.Zoom = 160
